# The "Client discount" line item (row 18) had a text note in the unit-price
# cell ("This client doesn't benefit from any discount") instead of a number,
# which made the dependent formulas (F18, F21, F23, F24) evaluate to #VALUE!.
# Replace it with a numeric 0 so the invoice totals calculate correctly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice")

$ws.Range("E18").Value = 0
